# Memory deallocation bug fix
# Adds bank-segment calculation rows to the "Dynamic" sheet and tidies up
# leftover selection/formatting state left over from earlier edits.

$wb = $excel.ActiveWorkbook

# --- Dynamic sheet: add the new "Bank Segment" calc block (rows 14-16) ---
$dyn = $wb.Worksheets.Item("Dynamic")

$dyn.Range("A14").Value = "Bank Segment"
$dyn.Range("B14").Value = 4

$dyn.Range("A15").Value = "Segments Per Bank"
$dyn.Range("B15").Formula = "=C4/F4"

$dyn.Range("A16").Value = "Segment Calc"
$dyn.Range("B16").Formula = "=(22-18) * B15 + B14"

# --- BANK61 sheet: selection had accidentally grown to the whole sheet; fix it ---
$bank61 = $wb.Worksheets.Item("BANK61")
$bank61.Range("B5:B6").Select()

# --- Golden sheet: clear the stray leftover bold formatting on A10 ---
$golden = $wb.Worksheets.Item("Golden")
$golden.Range("A10").Style = "Normal"

# --- Make "Dynamic" the active sheet again, with B16 selected ---
$dyn.Activate()
$dyn.Range("B16").Select()
